# Fruta / hortaliza, semanal
# Insert two new weekly rows (188-189) for "Vega Monumental Concepción - Frutilla",
# pushing the previously existing rows 188-192 down to 190-194.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 188, shifting existing rows 188:192 -> 190:194
$ws.Range("A188:A189").EntireRow.Insert()

# ---- New row 188 ----
$ws.Range("A188").Value = 11
$ws.Range("B188").Value = "Vega Monumental Concepción"
$ws.Range("C188").Value = "Bíobío"
$ws.Range("D188").Value = 44509
$ws.Range("E188").Value = 8
$ws.Range("F188").Value = "Fruta"
$ws.Range("G188").Value = 100101
$ws.Range("H188").Value = "Berries"
$ws.Range("I188").Value = 100112025
$ws.Range("J188").Value = "Frutilla"
$ws.Range("K188").Value = "Sin especificar"
$ws.Range("L188").Value = "Primera"
$ws.Range("M188").Value = 450
$ws.Range("N188").Value = 7500
$ws.Range("O188").Value = 8000
$ws.Range("P188").Value = 7722
$ws.Range("Q188").Value = "$/bandeja 7 kilos"
$ws.Range("R188").Value = "Provincia de Melipilla"
$ws.Range("S188").Value = 1103
$ws.Range("T188").Value = 7

# ---- New row 189 ----
$ws.Range("A189").Value = 11
$ws.Range("B189").Value = "Vega Monumental Concepción"
$ws.Range("C189").Value = "Bíobío"
$ws.Range("D189").Value = 44509
$ws.Range("E189").Value = 8
$ws.Range("F189").Value = "Fruta"
$ws.Range("G189").Value = 100101
$ws.Range("H189").Value = "Berries"
$ws.Range("I189").Value = 100112025
$ws.Range("J189").Value = "Frutilla"
$ws.Range("K189").Value = "Sin especificar"
$ws.Range("L189").Value = "Segunda"
$ws.Range("M189").Value = 250
$ws.Range("N189").Value = 6000
$ws.Range("O189").Value = 6500
$ws.Range("P189").Value = 6260
$ws.Range("Q189").Value = "$/bandeja 7 kilos"
$ws.Range("R189").Value = "Provincia de Melipilla"
$ws.Range("S189").Value = 894
$ws.Range("T189").Value = 7
